$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2523.2942
$ws.Range("I132").Value = 1859
$ws.Range("J132").Value = 3472.2856
$ws.Range("K132").Value = 5577
$ws.Range("L132").Value = 10416.8568
$ws.Range("M132").Value = -3047
$ws.Range("N132").Value = -15476.8568
$ws.Range("H137").Value = 1228.7561
$ws.Range("I137").Value = 802.8889
$ws.Range("J137").Value = 1348.5312
$ws.Range("K137").Value = 2408.6667
$ws.Range("L137").Value = 4045.5936
$ws.Range("M137").Value = 141.3332999999998
$ws.Range("N137").Value = -9145.5936
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6769.2075
$ws.Range("I32").Value = 3219.35
$ws.Range("J32").Value = 17691.846
$ws.Range("K32").Value = 3219.35
$ws.Range("L32").Value = 17691.846
$ws.Range("M32").Value = -2932.35
$ws.Range("N32").Value = -18265.846
$ws.Range("H45").Value = 1185.5714
$ws.Range("I45").Value = 930.8461
$ws.Range("J45").Value = 1599.5
$ws.Range("K45").Value = 930.8461
$ws.Range("L45").Value = 1599.5
$ws.Range("M45").Value = -553.8461
$ws.Range("N45").Value = -2353.5
$ws.Range("H61").Value = 1698.5
$ws.Range("I61").Value = 1139.3684
$ws.Range("J61").Value = 3216.1428
$ws.Range("K61").Value = 1139.3684
$ws.Range("L61").Value = 3216.1428
$ws.Range("M61").Value = -927.3684000000001
$ws.Range("N61").Value = -3640.1428
$ws.Range("H74").Value = 1107.4286
$ws.Range("I74").Value = 709.4666999999999
$ws.Range("J74").Value = 2102.3333
$ws.Range("K74").Value = 709.4666999999999
$ws.Range("L74").Value = 2102.3333
$ws.Range("M74").Value = 164.5333000000001
$ws.Range("N74").Value = -3850.3333
$ws.Range("H77").Value = 1107.4286
$ws.Range("I77").Value = 709.4666999999999
$ws.Range("J77").Value = 2102.3333
$ws.Range("K77").Value = 3547.3335
$ws.Range("L77").Value = 10511.6665
$ws.Range("M77").Value = 820.6665000000003
$ws.Range("N77").Value = -19247.6665
$ws.Range("H110").Value = 2978.5715
$ws.Range("I110").Value = 1500
$ws.Range("J110").Value = 3225
$ws.Range("K110").Value = 1500
$ws.Range("L110").Value = 3225
$ws.Range("M110").Value = 545
$ws.Range("N110").Value = -7315
$ws.Range("H136").Value = 1698.5
$ws.Range("I136").Value = 1139.3684
$ws.Range("J136").Value = 3216.1428
$ws.Range("K136").Value = 3418.1052
$ws.Range("L136").Value = 9648.428400000001
$ws.Range("M136").Value = -868.1052
$ws.Range("N136").Value = -14748.4284
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1137.5
$ws.Range("I94").Value = 1025
$ws.Range("J94").Value = 1250
$ws.Range("K94").Value = 1025
$ws.Range("L94").Value = 1250
$ws.Range("M94").Value = -574
$ws.Range("N94").Value = -2152
$ws.Range("H99").Value = 1622.8334
$ws.Range("I99").Value = 1849.2354
$ws.Range("J99").Value = 1073
$ws.Range("K99").Value = 1849.2354
$ws.Range("L99").Value = 1073
$ws.Range("M99").Value = -351.2354
$ws.Range("N99").Value = -4069
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1838.8889
$ws.Range("I122").Value = 1864.2858
$ws.Range("J122").Value = 1750
$ws.Range("K122").Value = 5592.857400000001
$ws.Range("L122").Value = 5250
$ws.Range("M122").Value = -3142.857400000001
$ws.Range("N122").Value = -10150
$ws.Range("H134").Value = 45456770
$ws.Range("I134").Value = 55557296
$ws.Range("J134").Value = 4400
$ws.Range("K134").Value = 166671888
$ws.Range("L134").Value = 13200
$ws.Range("M134").Value = -166669353
$ws.Range("N134").Value = -18270
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 10330.2
$ws.Range("I2").Value = 53.333332
$ws.Range("J2").Value = 14734.571
$ws.Range("K2").Value = 319.999992
$ws.Range("L2").Value = 88407.42600000001
$ws.Range("M2").Value = -206.999992
$ws.Range("N2").Value = -88633.42600000001
$ws.Range("H4").Value = 3061371.2
$ws.Range("I4").Value = 4762112
$ws.Range("J4").Value = 85075.164
$ws.Range("K4").Value = 14286336
$ws.Range("L4").Value = 255225.492
$ws.Range("M4").Value = -14286224
$ws.Range("N4").Value = -255449.492
$ws.Range("H6").Value = 168
$ws.Range("I6").Value = 107.5
$ws.Range("J6").Value = 248.66667
$ws.Range("K6").Value = 322.5
$ws.Range("L6").Value = 746.00001
$ws.Range("M6").Value = -209.5
$ws.Range("N6").Value = -972.00001
$ws.Range("H40").Value = 109.07692
$ws.Range("I40").Value = 109.07692
$ws.Range("K40").Value = 436.30768
$ws.Range("M40").Value = -367.30768
$ws.Range("H108").Value = 866.25
$ws.Range("I108").Value = 321.66666
$ws.Range("J108").Value = 2500
$ws.Range("K108").Value = 964.9999799999999
$ws.Range("L108").Value = 7500
$ws.Range("M108").Value = 1915.00002
$ws.Range("N108").Value = -13260
$ws.Range("H109").Value = 1668326.5
$ws.Range("I109").Value = 934.1429000000001
$ws.Range("J109").Value = 5558908.5
$ws.Range("K109").Value = 2802.4287
$ws.Range("L109").Value = 16676725.5
$ws.Range("M109").Value = -1762.4287
$ws.Range("N109").Value = -16678805.5
$ws.Range("H129").Value = 1150.7646
$ws.Range("I129").Value = 771.7
$ws.Range("J129").Value = 1692.2858
$ws.Range("K129").Value = 2315.1
$ws.Range("L129").Value = 5076.857400000001
$ws.Range("M129").Value = 2684.9
$ws.Range("N129").Value = -15076.8574
$ws.Range("H131").Value = 8929380
$ws.Range("J131").Value = 1113.9395
$ws.Range("L131").Value = 3341.8185
$ws.Range("N131").Value = -13421.8185
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3622.3076
$ws.Range("I80").Value = 2848.75
$ws.Range("J80").Value = 4860
$ws.Range("K80").Value = 2848.75
$ws.Range("L80").Value = 4860
$ws.Range("M80").Value = -1850.75
$ws.Range("N80").Value = -6856
$ws.Range("H83").Value = 3622.3076
$ws.Range("I83").Value = 2848.75
$ws.Range("J83").Value = 4860
$ws.Range("K83").Value = 14243.75
$ws.Range("L83").Value = 24300
$ws.Range("M83").Value = -9251.75
$ws.Range("N83").Value = -34284
$ws.Range("H113").Value = 9071.429
$ws.Range("I113").Value = 21460
$ws.Range("J113").Value = 2188.889
$ws.Range("K113").Value = 21460
$ws.Range("L113").Value = 2188.889
$ws.Range("M113").Value = -19290
$ws.Range("N113").Value = -6528.889
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 19124.875
$ws.Range("I81").Value = 100000
$ws.Range("J81").Value = 7571.2856
$ws.Range("K81").Value = 200000
$ws.Range("L81").Value = 15142.5712
$ws.Range("M81").Value = -198939
$ws.Range("N81").Value = -17264.5712
$ws.Range("H84").Value = 19124.875
$ws.Range("I84").Value = 100000
$ws.Range("J84").Value = 7571.2856
$ws.Range("K84").Value = 1000000
$ws.Range("L84").Value = 75712.856
$ws.Range("M84").Value = -994696
$ws.Range("N84").Value = -86320.856
$ws.Range("H126").Value = 1700
$ws.Range("I126").Value = 1400
$ws.Range("J126").Value = 2200
$ws.Range("K126").Value = 4200
$ws.Range("L126").Value = 6600
$ws.Range("M126").Value = -1730
$ws.Range("N126").Value = -11540
$ws.Range("H132").Value = 2018.6061
$ws.Range("I132").Value = 1255.6
$ws.Range("J132").Value = 3192.4614
$ws.Range("K132").Value = 3766.8
$ws.Range("L132").Value = 9577.3842
$ws.Range("M132").Value = -1236.8
$ws.Range("N132").Value = -14637.3842
